$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Split column A out of the combined A:B width group (same width, just its own <col>) ---
$ws.Columns.Item(1).ColumnWidth = 29.86

# --- Drop the old rows 10-22 block so it can be rebuilt with the new (longer) layout ---
$ws.Range("A10:C22").EntireRow.Delete()

# Row 9 (Semestre ideal:) already carries the three column styles we need to reuse:
# A -> bold label style, B -> wrapped value style, C -> wrapped red-value style.

# --- Row 10 ---
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10,1).Value = 'Objetivos:'
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10,2).Value = 'Fornecer os conhecimentos teóricos e práticos sobre estatística aplicada, materialografia e análises térmicas de materiais.'
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10,3).Value = 'Fornecer os conhecimentos teóricos e práticos sobre estatística aplicada, materialografia e análises térmicas de materiais.'
$ws.Rows.Item(10).RowHeight = 60

# --- Row 11 ---
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(11,1).Value = 'Objectives:'
$ws.Rows.Item(11).RowHeight = 60

# --- Row 12 ---
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(12,1).Value = 'Docentes responsáveis:'

# --- Row 13 ---
$ws.Range("B9").Copy()
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(13,2).Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(13,3).Value = '519033 - Carlos Yujiro Shigue'

# --- Row 14 ---
$ws.Range("B9").Copy()
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(14,2).Value = '5840963 - Daniela Camargo Vernilli'
$ws.Range("C9").Copy()
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(14,3).Value = '5840963 - Daniela Camargo Vernilli'

# --- Row 15 ---
$ws.Range("B9").Copy()
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,2).Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("C9").Copy()
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,3).Value = '6495737 - Durval Rodrigues Junior'

# --- Row 16 ---
$ws.Range("B9").Copy()
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16,2).Value = '984972 - Hugo Ricardo Zschommler Sandim'
$ws.Range("C9").Copy()
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16,3).Value = '984972 - Hugo Ricardo Zschommler Sandim'

# --- Row 17 ---
$ws.Range("A9").Copy()
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17,1).Value = 'Programa resumido:'
$ws.Range("B9").Copy()
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17,2).Value = 'Técnicas de Materialografia. Calorimetria e análises térmicas de materiais.'
$ws.Range("C9").Copy()
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17,3).Value = 'Técnicas de Materialografia. Calorimetria e análises térmicas de materiais.'
$ws.Rows.Item(17).RowHeight = 60

# --- Row 18 ---
$ws.Range("A9").Copy()
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(18,1).Value = 'Short syllabus:'
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19 ---
$ws.Range("A9").Copy()
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19,1).Value = 'Programa:'
$ws.Range("B9").Copy()
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19,2).Value = '1. MATERIALOGRAFIA: Ensaio macrográfico ou macrografia; ensaio micrográfico ou micrografia. CORPO DE PROVA OU AMOSTRA. CORTE: discos de corte. Procedimento para o corte. EMBUTIMENTO: Preparação de corpo de prova: corpo de prova embutido a quente e a frio. Corpo de prova não embutido. LIXAMENTO: tipos de lixa; procedimento para o lixamento. POLIMENTO: processo mecânico; cuidados a serem observados no polimento. Processo semiautomático; processo eletrolítico; processo mecânico eletrolítico; polimento químico. Escolha do tipo de polimento. Procedimento para o polimento. ATAQUE QUÍMICO: princípio; métodos para obtenção de contraste. MICROSCOPIA ÓPTICA: Iluminação campo escuro; luz polarizada; contraste de fase; interferência diferencial. Partes de um microscópio óptico de reflexão; elementos mecânicos; elementos ópticos; iluminador; acessórios. Princípio da formação da imagem. Microscópio óptico de reflexão.
2. CALORIMETRIA E ANÁLISES TÉRMICAS: Fundamentos termodinâmicos da calorimetria e análises térmicas. Princípios de calorimetria e tipos de calorímetros. Análise térmica diferencial (DTA) e calorimetria exploratória diferencial (DSC): princípios de DTA e DSC; tipos de equipamentos: DSC de compensação de energia e DSC de fluxo de calor. Aplicações de DTA e DSC. Equipamento: cadinhos de DTA; cadinhos de DSC. Cálculo de entalpia; linha base e cálculo de calor específico. Determinação de transição de fases. Determinação do diagrama de fases de ligas binárias por DTA/DSC. Cálculos cinéticos de cristalização, transições de fases e reações de polimerização. Termogravimetria (TGA): definição; aplicações da TGA. Equipamento: forno; programador de temperatura; termopar; balança; tipos de cadinho. Avaliação de estabilidade térmica e estudos de envelhecimento de polímeros. Técnicas de análises térmicas acopladas a análise de gases evolvidos por espectrometria de massa (TGA-MS) e FTIR (TGA-FTIR).'
$ws.Range("C9").Copy()
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19,3).Value = '1. MATERIALOGRAFIA: Ensaio macrográfico ou macrografia; ensaio micrográfico ou micrografia. CORPO DE PROVA OU AMOSTRA. CORTE: discos de corte. Procedimento para o corte. EMBUTIMENTO: Preparação de corpo de prova: corpo de prova embutido a quente e a frio. Corpo de prova não embutido. LIXAMENTO: tipos de lixa; procedimento para o lixamento. POLIMENTO: processo mecânico; cuidados a serem observados no polimento. Processo semiautomático; processo eletrolítico; processo mecânico eletrolítico; polimento químico. Escolha do tipo de polimento. Procedimento para o polimento. ATAQUE QUÍMICO: princípio; métodos para obtenção de contraste. MICROSCOPIA ÓPTICA: Iluminação campo escuro; luz polarizada; contraste de fase; interferência diferencial. Partes de um microscópio óptico de reflexão; elementos mecânicos; elementos ópticos; iluminador; acessórios. Princípio da formação da imagem. Microscópio óptico de reflexão.
2. CALORIMETRIA E ANÁLISES TÉRMICAS: Fundamentos termodinâmicos da calorimetria e análises térmicas. Princípios de calorimetria e tipos de calorímetros. Análise térmica diferencial (DTA) e calorimetria exploratória diferencial (DSC): princípios de DTA e DSC; tipos de equipamentos: DSC de compensação de energia e DSC de fluxo de calor. Aplicações de DTA e DSC. Equipamento: cadinhos de DTA; cadinhos de DSC. Cálculo de entalpia; linha base e cálculo de calor específico. Determinação de transição de fases. Determinação do diagrama de fases de ligas binárias por DTA/DSC. Cálculos cinéticos de cristalização, transições de fases e reações de polimerização. Termogravimetria (TGA): definição; aplicações da TGA. Equipamento: forno; programador de temperatura; termopar; balança; tipos de cadinho. Avaliação de estabilidade térmica e estudos de envelhecimento de polímeros. Técnicas de análises térmicas acopladas a análise de gases evolvidos por espectrometria de massa (TGA-MS) e FTIR (TGA-FTIR).'
$ws.Rows.Item(19).RowHeight = 120

# --- Row 20 ---
$ws.Range("A9").Copy()
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20,1).Value = 'Syllabus:'
$ws.Rows.Item(20).RowHeight = 120

# --- Row 21 ---
$ws.Range("A9").Copy()
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(21,1).Value = 'Avaliação:'

# --- Row 22 ---
$ws.Range("A9").Copy()
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22,1).Value = 'Método:'
$ws.Range("B9").Copy()
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22,2).Value = 'Aulas expositivas complementadas com experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento e de estudo de casos.'
$ws.Range("C9").Copy()
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22,3).Value = 'Aulas expositivas complementadas com experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento e de estudo de casos.'
$ws.Rows.Item(22).RowHeight = 60

# --- Row 23 ---
$ws.Range("A9").Copy()
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23,1).Value = 'Critério:'
$ws.Range("B9").Copy()
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23,2).Value = 'Média aritmética das notas obtidas nos relatórios e trabalhos. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0.'
$ws.Range("C9").Copy()
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23,3).Value = 'Média aritmética das notas obtidas nos relatórios e trabalhos. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0.'
$ws.Rows.Item(23).RowHeight = 60

# --- Row 24 ---
$ws.Range("A9").Copy()
$ws.Range("A24").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(24,1).Value = 'Norma de recuperação:'
$ws.Range("B9").Copy()
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(24,2).Value = 'Devido às características práticas da disciplina, não será oferecida recuperação.'
$ws.Range("C9").Copy()
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(24,3).Value = 'Devido às características práticas da disciplina, não será oferecida recuperação.'
$ws.Rows.Item(24).RowHeight = 60

# --- Row 25 ---
$ws.Range("A9").Copy()
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25,1).Value = 'Bibliografia:'
$ws.Range("B9").Copy()
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25,2).Value = 'COLPAERT; HUBERTUS. Metalografia dos produtos siderúrgicos comuns, 3ª Edição, Editora Edgard Blücher Ltda, São Paulo – 1974.COUTINHO, TELMO DE AZEVEDO. Metalografia de Não-Ferrosos, Editora Edgard Blücher Ltda, São Paulo – 1980.PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.AZEVEDO, A. D.; MOTHE, C. G. Análaise Térmica de Materiais. São Paulo: ARTLIBER, 2009.BROWN, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, New York: Wiley, 1999.HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.MULLER, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.SPEYER, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.REED-HILL, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. Nondestructive Characterization of Materials. Series. Plenum Press, New York. YACOBI, B.G.; HOLT, D.B.; KAZMERSKI, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994.'
$ws.Range("C9").Copy()
$ws.Range("C25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25,3).Value = 'COLPAERT; HUBERTUS. Metalografia dos produtos siderúrgicos comuns, 3ª Edição, Editora Edgard Blücher Ltda, São Paulo – 1974.COUTINHO, TELMO DE AZEVEDO. Metalografia de Não-Ferrosos, Editora Edgard Blücher Ltda, São Paulo – 1980.PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.AZEVEDO, A. D.; MOTHE, C. G. Análaise Térmica de Materiais. São Paulo: ARTLIBER, 2009.BROWN, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, New York: Wiley, 1999.HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.MULLER, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.SPEYER, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.REED-HILL, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. Nondestructive Characterization of Materials. Series. Plenum Press, New York. YACOBI, B.G.; HOLT, D.B.; KAZMERSKI, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994.'
$ws.Rows.Item(25).RowHeight = 120

# --- Row 26 ---
$ws.Range("A9").Copy()
$ws.Range("A26").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(26,1).Value = 'Requisitos:'

# --- Row 27 ---
$ws.Range("B9").Copy()
$ws.Range("B27").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27,2).Value = 'LOB1012 -  Estatística  (Requisito fraco)
'
$ws.Range("C9").Copy()
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27,3).Value = 'LOB1012 -  Estatística  (Requisito fraco)
'
$ws.Rows.Item(27).RowHeight = 30

$excel.CutCopyMode = 0

"done"